# cryptos.xlsx - scheduled data refresh (GitHub Actions)
#
# Rewrites the Coin / Link / Price / Volume(1h) columns on Sheet1 with the
# latest coinranking.com snapshot. Rows 29-31 also got re-ranked, so those
# three rows have their Coin (B) and Link (C) cells rewritten too, not just
# Price (D) and Volume(1h) (E).
#
# Price/Volume values are stored as plain TEXT, not numbers (prices use dots
# as thousands separators, e.g. "62.381.26", and volumes are space-padded
# percent strings, e.g. "  +2.56%  "). Assigning those strings straight to
# .Value would let Excel reinterpret the numeric-looking ones (e.g. "1.00",
# "5.40") as numbers and silently drop the trailing zero, so every write is
# prefixed with a literal leading apostrophe to force text, exactly like typing
# it into the grid. The apostrophe itself never becomes part of the stored
# value. Writing a leading apostrophe also nudges Excel to tag the cell with a
# "Text" number format, so .Style is reset to "Normal" right after to keep the
# cell's format exactly as it was (unstyled/general) before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$v = "'" + '62.381.26'
$ws.Range("D2").Value = $v
$ws.Range("D2").Style = "Normal"
$v = "'" + '  +2.56%  '
$ws.Range("E2").Value = $v
$ws.Range("E2").Style = "Normal"

# Row 3
$v = "'" + '2.424.70'
$ws.Range("D3").Value = $v
$ws.Range("D3").Style = "Normal"
$v = "'" + '  +3.22%  '
$ws.Range("E3").Value = $v
$ws.Range("E3").Style = "Normal"

# Row 4
$v = "'" + '  -0.02%  '
$ws.Range("E4").Value = $v
$ws.Range("E4").Style = "Normal"

# Row 5
$v = "'" + '556.45'
$ws.Range("D5").Value = $v
$ws.Range("D5").Style = "Normal"
$v = "'" + '  +2.26%  '
$ws.Range("E5").Value = $v
$ws.Range("E5").Style = "Normal"

# Row 6
$v = "'" + '143.50'
$ws.Range("D6").Value = $v
$ws.Range("D6").Style = "Normal"
$v = "'" + '  +4.99%  '
$ws.Range("E6").Value = $v
$ws.Range("E6").Style = "Normal"

# Row 7
$v = "'" + '  -0.05%  '
$ws.Range("E7").Value = $v
$ws.Range("E7").Style = "Normal"

# Row 8
$v = "'" + '  +1.77%  '
$ws.Range("E8").Value = $v
$ws.Range("E8").Style = "Normal"

# Row 9
$v = "'" + '2.424.73'
$ws.Range("D9").Value = $v
$ws.Range("D9").Style = "Normal"
$v = "'" + '  +3.31%  '
$ws.Range("E9").Value = $v
$ws.Range("E9").Style = "Normal"

# Row 10
$v = "'" + '  +4.96%  '
$ws.Range("E10").Value = $v
$ws.Range("E10").Style = "Normal"

# Row 11
$v = "'" + '  -0.49%  '
$ws.Range("E11").Value = $v
$ws.Range("E11").Style = "Normal"

# Row 12
$v = "'" + '5.40'
$ws.Range("D12").Value = $v
$ws.Range("D12").Style = "Normal"
$v = "'" + '  +1.66%  '
$ws.Range("E12").Value = $v
$ws.Range("E12").Style = "Normal"

# Row 13
$v = "'" + '  +2.61%  '
$ws.Range("E13").Value = $v
$ws.Range("E13").Style = "Normal"

# Row 14
$v = "'" + '26.33'
$ws.Range("D14").Value = $v
$ws.Range("D14").Style = "Normal"
$v = "'" + '  +6.62%  '
$ws.Range("E14").Value = $v
$ws.Range("E14").Style = "Normal"

# Row 15
$v = "'" + '  +9.62%  '
$ws.Range("E15").Value = $v
$ws.Range("E15").Style = "Normal"

# Row 16
$v = "'" + '2.864.29'
$ws.Range("D16").Value = $v
$ws.Range("D16").Style = "Normal"
$v = "'" + '  +3.23%  '
$ws.Range("E16").Value = $v
$ws.Range("E16").Style = "Normal"

# Row 17
$v = "'" + '62.345.47'
$ws.Range("D17").Value = $v
$ws.Range("D17").Style = "Normal"
$v = "'" + '  +2.52%  '
$ws.Range("E17").Value = $v
$ws.Range("E17").Style = "Normal"

# Row 18
$v = "'" + '2.425.74'
$ws.Range("D18").Value = $v
$ws.Range("D18").Style = "Normal"
$v = "'" + '  +3.32%  '
$ws.Range("E18").Value = $v
$ws.Range("E18").Style = "Normal"

# Row 19
$v = "'" + '11.09'
$ws.Range("D19").Value = $v
$ws.Range("D19").Style = "Normal"
$v = "'" + '  +4.17%  '
$ws.Range("E19").Value = $v
$ws.Range("E19").Style = "Normal"

# Row 20
$v = "'" + '324.71'
$ws.Range("D20").Value = $v
$ws.Range("D20").Style = "Normal"
$v = "'" + '  +1.82%  '
$ws.Range("E20").Value = $v
$ws.Range("E20").Style = "Normal"

# Row 21
$v = "'" + '  +1.70%  '
$ws.Range("E21").Value = $v
$ws.Range("E21").Style = "Normal"

# Row 22
$v = "'" + '6.74'
$ws.Range("D22").Value = $v
$ws.Range("D22").Style = "Normal"
$v = "'" + '  +3.00%  '
$ws.Range("E22").Value = $v
$ws.Range("E22").Style = "Normal"

# Row 23
$v = "'" + '  +0.34%  '
$ws.Range("E23").Value = $v
$ws.Range("E23").Style = "Normal"

# Row 24
$v = "'" + '1.78'
$ws.Range("D24").Value = $v
$ws.Range("D24").Style = "Normal"
$v = "'" + '  +5.99%  '
$ws.Range("E24").Value = $v
$ws.Range("E24").Style = "Normal"

# Row 25
$v = "'" + '64.99'
$ws.Range("D25").Value = $v
$ws.Range("D25").Style = "Normal"
$v = "'" + '  +2.64%  '
$ws.Range("E25").Value = $v
$ws.Range("E25").Style = "Normal"

# Row 26
$v = "'" + '  +9.18%  '
$ws.Range("E26").Value = $v
$ws.Range("E26").Style = "Normal"

# Row 27
$v = "'" + '574.00'
$ws.Range("D27").Value = $v
$ws.Range("D27").Style = "Normal"
$v = "'" + '  +14.94%  '
$ws.Range("E27").Value = $v
$ws.Range("E27").Style = "Normal"

# Row 28
$v = "'" + '2.547.65'
$ws.Range("D28").Value = $v
$ws.Range("D28").Style = "Normal"
$v = "'" + '  +3.29%  '
$ws.Range("E28").Value = $v
$ws.Range("E28").Style = "Normal"

# Row 29
$v = "'" + 'Binance-PegBSC-USD'
$ws.Range("B29").Value = $v
$ws.Range("B29").Style = "Normal"
$v = "'" + 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("C29").Value = $v
$ws.Range("C29").Style = "Normal"
$v = "'" + '1.00'
$ws.Range("D29").Value = $v
$ws.Range("D29").Style = "Normal"
$v = "'" + '  -0.13%  '
$ws.Range("E29").Value = $v
$ws.Range("E29").Style = "Normal"

# Row 30
$v = "'" + 'InternetComputer(DFINITY)'
$ws.Range("B30").Value = $v
$ws.Range("B30").Style = "Normal"
$v = "'" + 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C30").Value = $v
$ws.Range("C30").Style = "Normal"
$v = "'" + '8.40'
$ws.Range("D30").Value = $v
$ws.Range("D30").Style = "Normal"
$v = "'" + '  +5.59%  '
$ws.Range("E30").Value = $v
$ws.Range("E30").Style = "Normal"

# Row 31
$v = "'" + 'PEPE'
$ws.Range("B31").Value = $v
$ws.Range("B31").Style = "Normal"
$v = "'" + 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C31").Value = $v
$ws.Range("C31").Style = "Normal"
$v = "'" + '0.0₃0941'
$ws.Range("D31").Value = $v
$ws.Range("D31").Style = "Normal"
$v = "'" + '  +9.74%  '
$ws.Range("E31").Value = $v
$ws.Range("E31").Style = "Normal"

# Row 32
$v = "'" + '1.45'
$ws.Range("D32").Value = $v
$ws.Range("D32").Style = "Normal"
$v = "'" + '  +6.02%  '
$ws.Range("E32").Value = $v
$ws.Range("E32").Style = "Normal"

# Row 33
$v = "'" + '0.148'
$ws.Range("D33").Value = $v
$ws.Range("D33").Style = "Normal"
$v = "'" + '  +2.10%  '
$ws.Range("E33").Value = $v
$ws.Range("E33").Style = "Normal"

# Row 34
$v = "'" + '  +4.21%  '
$ws.Range("E34").Value = $v
$ws.Range("E34").Style = "Normal"

# Row 35
$v = "'" + '1.57'
$ws.Range("D35").Value = $v
$ws.Range("D35").Style = "Normal"
$v = "'" + '  +5.47%  '
$ws.Range("E35").Value = $v
$ws.Range("E35").Style = "Normal"

# Row 36
$v = "'" + '5.73'
$ws.Range("D36").Value = $v
$ws.Range("D36").Style = "Normal"
$v = "'" + '  +9.36%  '
$ws.Range("E36").Value = $v
$ws.Range("E36").Style = "Normal"

# Row 37
$v = "'" + '4.85'
$ws.Range("D37").Value = $v
$ws.Range("D37").Style = "Normal"
$v = "'" + '  +5.58%  '
$ws.Range("E37").Value = $v
$ws.Range("E37").Style = "Normal"

# Row 38
$v = "'" + '  -0.05%  '
$ws.Range("E38").Value = $v
$ws.Range("E38").Style = "Normal"

# Row 39
$v = "'" + '  +2.55%  '
$ws.Range("E39").Value = $v
$ws.Range("E39").Style = "Normal"

# Row 40
$v = "'" + '1.88'
$ws.Range("D40").Value = $v
$ws.Range("D40").Style = "Normal"
$v = "'" + '  +2.82%  '
$ws.Range("E40").Value = $v
$ws.Range("E40").Style = "Normal"

# Row 41
$v = "'" + '18.78'
$ws.Range("D41").Value = $v
$ws.Range("D41").Style = "Normal"
$v = "'" + '  +1.73%  '
$ws.Range("E41").Value = $v
$ws.Range("E41").Style = "Normal"

# Row 42
$v = "'" + '150.46'
$ws.Range("D42").Value = $v
$ws.Range("D42").Style = "Normal"
$v = "'" + '  +5.05%  '
$ws.Range("E42").Value = $v
$ws.Range("E42").Style = "Normal"

# Row 44
$v = "'" + '41.70'
$ws.Range("D44").Value = $v
$ws.Range("D44").Style = "Normal"
$v = "'" + '  +2.72%  '
$ws.Range("E44").Value = $v
$ws.Range("E44").Style = "Normal"

# Row 45
$v = "'" + '2.34'
$ws.Range("D45").Value = $v
$ws.Range("D45").Style = "Normal"
$v = "'" + '  +14.91%  '
$ws.Range("E45").Value = $v
$ws.Range("E45").Style = "Normal"

# Row 46
$v = "'" + '151.08'
$ws.Range("D46").Value = $v
$ws.Range("D46").Style = "Normal"
$v = "'" + '  +5.65%  '
$ws.Range("E46").Value = $v
$ws.Range("E46").Style = "Normal"

# Row 47
$v = "'" + '3.64'
$ws.Range("D47").Value = $v
$ws.Range("D47").Style = "Normal"
$v = "'" + '  +2.39%  '
$ws.Range("E47").Value = $v
$ws.Range("E47").Style = "Normal"

# Row 48
$v = "'" + '0.0543'
$ws.Range("D48").Value = $v
$ws.Range("D48").Style = "Normal"
$v = "'" + '  +5.11%  '
$ws.Range("E48").Value = $v
$ws.Range("E48").Style = "Normal"

# Row 49
$v = "'" + '20.45'
$ws.Range("D49").Value = $v
$ws.Range("D49").Style = "Normal"
$v = "'" + '  +7.37%  '
$ws.Range("E49").Value = $v
$ws.Range("E49").Style = "Normal"

# Row 50
$v = "'" + '0.590'
$ws.Range("D50").Value = $v
$ws.Range("D50").Style = "Normal"
$v = "'" + '  +3.98%  '
$ws.Range("E50").Value = $v
$ws.Range("E50").Style = "Normal"

# Row 51
$v = "'" + '0.0917'
$ws.Range("D51").Value = $v
$ws.Range("D51").Style = "Normal"
$v = "'" + '  +1.85%  '
$ws.Range("E51").Value = $v
$ws.Range("E51").Style = "Normal"

